# Start to prepare and cleanse the raw "upah-minimum" dataset:
# the second column was a generic "value" header - replace it with the
# actual year the minimum-wage figures belong to (2023), stored as a
# real number instead of the placeholder text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 2023
